$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

$timestamps = @(
    "2021-10-05 13:42:22.727704",
    "2021-10-05 13:42:22.727716",
    "2021-10-05 13:42:22.727720",
    "2021-10-05 13:42:22.727722",
    "2021-10-05 13:42:22.727725",
    "2021-10-05 13:42:22.727728",
    "2021-10-05 13:42:22.727730",
    "2021-10-05 13:42:22.727733",
    "2021-10-05 13:42:22.727736",
    "2021-10-05 13:42:22.727738",
    "2021-10-05 13:42:22.727741",
    "2021-10-05 13:42:22.727743",
    "2021-10-05 13:42:22.727746",
    "2021-10-05 13:42:22.727748",
    "2021-10-05 13:42:22.727751",
    "2021-10-05 13:42:22.727753",
    "2021-10-05 13:42:22.727756",
    "2021-10-05 13:42:22.727759",
    "2021-10-05 13:42:22.727761",
    "2021-10-05 13:42:22.727764",
    "2021-10-05 13:42:22.727766",
    "2021-10-05 13:42:22.727769",
    "2021-10-05 13:42:22.727771",
    "2021-10-05 13:42:22.727774",
    "2021-10-05 13:42:22.727777",
    "2021-10-05 13:42:22.727779",
    "2021-10-05 13:42:22.727782",
    "2021-10-05 13:42:22.727784",
    "2021-10-05 13:42:22.727787",
    "2021-10-05 13:42:22.727791",
    "2021-10-05 13:42:22.727794",
    "2021-10-05 13:42:22.727796",
    "2021-10-05 13:42:22.727799",
    "2021-10-05 13:42:22.727802",
    "2021-10-05 13:42:22.727804"
)

for ($i = 0; $i -lt $timestamps.Length; $i++) {
    $row = $i + 2
    $cell = $ws.Cells.Item($row, 6)
    $cell.Value = $timestamps[$i]
}
